$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "months" -> "maturities"
$ws.Range("B1").Value = "maturities"

# Column B: numeric month counts -> text "<n>M" maturity labels.
# (This also clears the B20/B28/B30 formulas, since the cells become
# plain text strings instead of numbers.)
$maturities = [ordered]@{
    2  = "1M"
    3  = "2M"
    4  = "3M"
    5  = "4M"
    6  = "5M"
    7  = "6M"
    8  = "7M"
    9  = "8M"
    10 = "9M"
    11 = "10M"
    12 = "11M"
    13 = "12M"
    14 = "15M"
    15 = "18M"
    16 = "21M"
    17 = "24M"
    18 = "36M"
    19 = "48M"
    20 = "60M"
    21 = "72M"
    22 = "84M"
    23 = "96M"
    24 = "108M"
    25 = "120M"
    26 = "132M"
    27 = "144M"
    28 = "180M"
    29 = "216M"
    30 = "300M"
    31 = "360M"
    32 = "480M"
    33 = "600M"
    34 = "720M"
}

foreach ($row in $maturities.Keys) {
    $ws.Cells.Item($row, 2).Value = $maturities[$row]
}

# Column C (rows 7-34) picks up the same style already used by C2:C6
# (Menlo 14pt) so every data row is consistently styled. Copy the
# formatting from C2 rather than re-assigning font properties one by
# one, so the existing style is reused instead of minting new ones.
$ws.Range("C2").Copy()
$ws.Range("C7:C34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Column A width adjusts to fit the new (shorter) content.
$ws.Columns.Item(1).AutoFit()

# Selection moves to the (now populated) column C.
$ws.Range("C:C").Select()
